# atualizei dados da bibi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13 (pushes existing rows 13+ down by one), adding
# a new daily entry for day 12 of 06/2025.
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 38843.63
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(13, 4).Value = 2025
$ws.Cells.Item(13, 5).Value = "06/2025"
